$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(4.453082120904227, 16.80072893514018, 31.66688530865395, 47.97825748599409, 14.99599398427064)
    3  = @(4.256110385768206, 15.81414347791592, 30.13649976548376, 45.07585997036804, 14.56739984798546)
    4  = @(4.129921971931085, 15.18345797957188, 29.16934848132695, 43.21230408299174, 14.30402308097856)
    5  = @(4.077210539987675, 14.92044224378863, 28.76891537651168, 42.43304870019231, 14.19681494693909)
    6  = @(4.068381018414702, 14.87641457860606, 28.70206186205785, 42.30247733761892, 14.17902561583316)
    7  = @(4.12921626029833,  15.17993479371134, 29.16397279584879, 43.20187414878651, 14.30257651222955)
    8  = @(4.386272154184288, 16.46586561567113, 31.14524790463043, 46.99476461840644, 14.84836263646835)
    9  = @(4.847626781533299, 18.9696309584411,  34.79019700568817, 53.76441982954865, 15.90964749827232)
    10 = @(5.159375015272057, 20.70052647128698, 37.29712272872885, 58.30908890538312, 16.67484774318619)
    11 = @(5.295119806458881, 21.4463021710597,  38.39657939506924, 60.280539084373,   17.0181105394761)
    12 = @(5.345639535608059, 21.72277565408117, 38.8067745188044,  61.01311911472015, 17.14727530363974)
    13 = @(5.334798663661791, 21.66349534754373, 38.71870855538845, 60.85596817346873, 17.11949578527255)
    14 = @(5.299293861980868, 21.46916648747221, 38.43045085574378, 60.34108956117134, 17.02875414856349)
    15 = @(5.277430761849376, 21.34936291773653, 38.25307759903464, 60.02388868579958, 16.97306174473766)
    16 = @(5.150380618130093, 20.65095500525957, 37.22442356079723, 58.17830377378787, 16.65230639989398)
    17 = @(5.070875926456229, 20.21187991088265, 36.58269053228319, 57.02138715638991, 16.45420217812308)
    18 = @(5.024575959885103, 19.95541536592154, 36.20973971098359, 56.34693520985561, 16.33981001259283)
    19 = @(5.008801897977285, 19.86790610021647, 36.08281351125076, 56.11703346922731, 16.30100592805423)
    20 = @(5.079398517153265, 20.25902557245059, 36.65140384155783, 57.14547825128371, 16.47533801822654)
    21 = @(5.309746544556227, 21.52640628333934, 38.51528771957194, 60.49270200548483, 17.05543041105215)
    22 = @(5.455137147564244, 22.32013858859543, 39.69753703054364, 62.59886043891652, 17.42971607441617)
    23 = @(5.378014062561181, 21.89965602861503, 39.06990661277875, 61.48225963977133, 17.23043429362994)
    24 = @(5.075547295935126, 20.23772358888641, 36.62035102466099, 57.08940569978326, 16.46578405071606)
    25 = @(4.72752760016914,  18.2949846601868,  33.83266512629793, 52.0071702442389,  15.62447975441803)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("D$row").Value = $values[0]
    $ws.Range("E$row").Value = $values[1]
    $ws.Range("F$row").Value = $values[2]
    $ws.Range("G$row").Value = $values[3]
    $ws.Range("H$row").Value = $values[4]
}
